$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.030.65"
$ws.Range("E2").Value = "  +1.52%  "

# Row 3
$ws.Range("D3").Value = "3.928.93"
$ws.Range("E3").Value = "  +0.56%  "

# Row 4
$ws.Range("E4").Value = "  +0.13%  "

# Row 5
$ws.Range("D5").Value = "'486.65"
$ws.Range("E5").Value = "  +0.40%  "

# Row 6
$ws.Range("D6").Value = "'147.57"
$ws.Range("E6").Value = "  +1.75%  "

# Row 7
$ws.Range("E7").Value = "  -0.28%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("D9").Value = "'0.739"
$ws.Range("E9").Value = "  +1.74%  "

# Row 10
$ws.Range("D10").Value = "'0.178"
$ws.Range("E10").Value = "  +6.53%  "

# Row 11
$ws.Range("D11").Value = "'0.0000348"
$ws.Range("E11").Value = "  -1.90%  "

# Row 12
$ws.Range("D12").Value = "'43.00"
$ws.Range("E12").Value = "  +1.22%  "

# Row 13
$ws.Range("E13").Value = "  -1.31%  "

# Row 14
$ws.Range("D14").Value = "4.557.85"
$ws.Range("E14").Value = "  +0.51%  "

# Row 15
$ws.Range("D15").Value = "3.940.14"
$ws.Range("E15").Value = "  +0.72%  "

# Row 16
$ws.Range("D16").Value = "'14.25"
$ws.Range("E16").Value = "  -3.32%  "

# Row 17
$ws.Range("E17").Value = "  -0.70%  "

# Row 18
$ws.Range("D18").Value = "'19.99"
$ws.Range("E18").Value = "  +0.97%  "

# Row 19
$ws.Range("E19").Value = "  +1.47%  "

# Row 20
$ws.Range("D20").Value = "69.082.67"
$ws.Range("E20").Value = "  +1.52%  "

# Row 21
$ws.Range("D21").Value = "'439.29"
$ws.Range("E21").Value = "  -1.74%  "

# Row 22
$ws.Range("D22").Value = "'3.49"
$ws.Range("E22").Value = "  +3.85%  "

# Row 23
$ws.Range("E23").Value = "  +0.33%  "

# Row 24
$ws.Range("D24").Value = "'89.66"
$ws.Range("E24").Value = "  +0.66%  "

# Row 25
$ws.Range("D25").Value = "'12.07"
$ws.Range("E25").Value = "  +9.76%  "

# Row 26
$ws.Range("D26").Value = "'3.71"
$ws.Range("E26").Value = "  +3.15%  "

# Row 27
$ws.Range("D27").Value = "'11.13"
$ws.Range("E27").Value = "  -4.82%  "

# Row 28
$ws.Range("D28").Value = "'37.18"
$ws.Range("E28").Value = "  -3.76%  "

# Row 29
$ws.Range("D29").Value = "'5.66"
$ws.Range("E29").Value = "  -2.92%  "

# Row 30
$ws.Range("D30").Value = "'713.02"
$ws.Range("E30").Value = "  +3.96%  "

# Row 31
$ws.Range("E31").Value = "  +0.94%  "

# Row 32
$ws.Range("D32").Value = "'13.39"
$ws.Range("E32").Value = "  +0.31%  "

# Row 33
$ws.Range("D33").Value = "'2.89"
$ws.Range("E33").Value = "  +1.43%  "

# Row 34
$ws.Range("B34").Value = "TheGraph"
$ws.Range("C34").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D34").Value = "'0.468"
$ws.Range("E34").Value = "  +30.78%  "

# Row 35
$ws.Range("B35").Value = "PEPE"
$ws.Range("C35").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D35").Value = "0.0₃0919"
$ws.Range("E35").Value = "  -1.92%  "

# Row 36
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "'6.04"
$ws.Range("E36").Value = "  +4.72%  "

# Row 37
$ws.Range("D37").Value = "'40.99"
$ws.Range("E37").Value = "  -1.23%  "

# Row 38
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").Value = "'60.90"
$ws.Range("E38").Value = "  +3.82%  "

# Row 39
$ws.Range("D39").Value = "'0.149"
$ws.Range("E39").Value = "  -0.15%  "

# Row 40
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  +0.14%  "

# Row 41
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.12%  "

# Row 42
$ws.Range("D42").Value = "'2.97"
$ws.Range("E42").Value = "  +3.36%  "

# Row 43
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0488"
$ws.Range("E43").Value = "  +2.53%  "

# Row 44
$ws.Range("B44").Value = "ThetaToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D44").Value = "'3.08"
$ws.Range("E44").Value = "  +1.40%  "

# Row 45
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "'3.01"
$ws.Range("E45").Value = "  -0.04%  "

# Row 46
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₆0374"
$ws.Range("E46").Value = "  +12.61%  "

# Row 47
$ws.Range("D47").Value = "'3.42"
$ws.Range("E47").Value = "  +8.37%  "

# Row 48
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "'0.143"
$ws.Range("E48").Value = "  +0.98%  "

# Row 49
$ws.Range("E49").Value = "  -1.97%  "

# Row 50
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "'2.96"
$ws.Range("E50").Value = "  +4.77%  "

# Row 51
$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").Value = "'2.08"
$ws.Range("E51").Value = "  -1.71%  "
